$d = $word.ActiveDocument

# Locate the target paragraph: the one ending with the migration explanation
# text, which currently carries a stray <w:u w:val="single"/> on its paragraph
# mark and is the last paragraph in the document body.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*no contaba con una provincia.*") {
        $target = $para
    }
}

if ($null -eq $target) {
    throw "Could not locate target paragraph"
}

$xml = '<w:p w:rsidR="003A5C00" w:rsidRPr="003A5C00" w:rsidRDefault="003A5C00" w:rsidP="003A5C00"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:tab/><w:t xml:space="preserve">Al migrar los datos de los clientes de la tabla </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Maestra</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> hacia la tabla </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Clientes</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">, existe la nulidad en el campo </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Provincia</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>, por lo cual se carga por defecto con el atributo “</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="es-ES"/></w:rPr><w:t>Migrada</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>”. De esta forma, dejamos en claro, que este cliente fue introducido al sistema a través de la migración y no contaba con una provincia.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Ttulo"/><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:lastRenderedPageBreak/><w:t>Criterios para campos repetidos</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:tab/><w:t xml:space="preserve">Dentro de la tabla </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Maestra</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> se encontraron 4 clientes, los cuales tienen repetido el DNI, por lo cual, la estrategia optada por el grupo, fue la de a uno de esos repetidos multiplicarlo por -1. De esta manera, quedan los dos DNI con el valor original inalterado, pero podemos diferenciar de manera univoca a uno del otro.</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:tab/><w:t>Llegado el caso de que uno de estos casos particulares se quiera modificar, no se dejaría modificar el usuario sin antes haber modificado el DNI negativo, dado que la aplicación Desktop, valida que el DNI además de ser numérico, sea positivo.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>'

$target.Range.InsertXML($xml)

Write-Output "done"
